$wb = $excel.ActiveWorkbook

# --- Add "Day 24" as a new sheet at the very end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws24 = $wb.Worksheets.Add($null, $lastSheet)
$ws24.Name = "Day 24"

# Row 1 (headers) entered first
$ws24.Range("A1").Value = "student"
$ws24.Range("B1").Value = "class"

# Then column A filled top to bottom
$ws24.Range("A2").Value = "A"
$ws24.Range("A3").Value = "B"
$ws24.Range("A4").Value = "C"
$ws24.Range("A5").Value = "D"
$ws24.Range("A6").Value = "E"
$ws24.Range("A7").Value = "F"
$ws24.Range("A8").Value = "G"
$ws24.Range("A9").Value = "H"
$ws24.Range("A10").Value = "I"

# Then column B filled top to bottom
$ws24.Range("B2").Value = "Math"
$ws24.Range("B3").Value = "English"
$ws24.Range("B4").Value = "Math"
$ws24.Range("B5").Value = "Biology"
$ws24.Range("B6").Value = "Math"
$ws24.Range("B7").Value = "Computer"
$ws24.Range("B8").Value = "Math"
$ws24.Range("B9").Value = "Math"
$ws24.Range("B10").Value = "Math"

$ws24.Range("B11").Select() | Out-Null

# --- Add "Day 25" as a new sheet after "Day 24" ---
$ws25 = $wb.Worksheets.Add($null, $ws24)
$ws25.Name = "Day 25"

$ws25.Range("A1").Value = "order_number"
$ws25.Range("B1").Value = "customer_number"

$ws25.Range("A2").Value = 1
$ws25.Range("A3").Value = 2
$ws25.Range("A4").Value = 3
$ws25.Range("A5").Value = 4

$ws25.Range("B2").Value = 1
$ws25.Range("B3").Value = 2
$ws25.Range("B4").Value = 3
$ws25.Range("B5").Value = 3

$ws25.Columns.Item(1).AutoFit()
$ws25.Columns.Item(2).AutoFit()

$ws25.Range("F8").Select() | Out-Null
